# Add the new "2020" column (Q) to the consumer price index table, mirroring
# the formatting already applied to the 2019 column (P), then restore the
# author's final selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newYearValues = @{
    4  = 2020
    5  = 109.7221295941265
    6  = 108.44905375816947
    7  = 109.90982951756889
    8  = 108.40606487500015
    9  = 109.40161876466024
    10 = 107.71155656686271
    11 = 111.78921596090774
    12 = 111.39254046803097
    13 = 110.44919152842827
    14 = 106.89826464456031
}

foreach ($row in 4..14) {
    $srcCell = $ws.Cells.Item($row, 16)  # column P (2019)
    $dstCell = $ws.Cells.Item($row, 17)  # column Q (2020)

    # Copy P's style/number-format onto Q, then write the 2020 value so the
    # literal number isn't clobbered by the paste.
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null
    $dstCell.Value = $newYearValues[$row]
}

$excel.CutCopyMode = $false

# Restore the saved selection/active cell.
$ws.Range("N14").Select() | Out-Null
